$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "22.478.99"
$ws.Range("E2").Value = "  +0.46%  "

$ws.Range("D3").Value = "1.571.46"
$ws.Range("E3").Value = "  +0.25%  "

$ws.Range("E4").Value = "  +0.00%  "

$ws.Range("E5").Value = "  +0.02%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "291.32"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.17%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3712"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -1.50%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "49.95"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +1.89%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.3377"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.54%  "

$ws.Range("E10").Value = "  +0.76%  "

$ws.Range("E11").Value = "  -0.58%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "21.14"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.64%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.013"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.91%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.955"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.63%  "

$ws.Range("D16").Value = "1.570.75"
$ws.Range("E16").Value = "  +0.35%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001119"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.76%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "90.49"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.84%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06782"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.50%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.337"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +2.18%  "

$ws.Range("E22").Value = "  -1.01%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "12.21"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +2.33%  "

$ws.Range("D24").Value = "22.469.30"
$ws.Range("E24").Value = "  +0.40%  "

$ws.Range("E26").Value = "  -3.34%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "20.03"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.95%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "149.11"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.78%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.068"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.97%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "125.16"
$ws.Range("D30").Style = "Normal"

$ws.Range("D31").Value = "1.747.67"

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.068"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +7.54%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "6.201"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +2.39%  "

$ws.Range("E34").Value = "  -0.29%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "9.778"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -3.22%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.08343"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.35%  "

$ws.Range("E37").Value = "  -4.70%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02478"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.45%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.2302"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.47%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.06541"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +1.39%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.433"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.20%  "

$ws.Range("E42").Value = "  +0.10%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.6201"
$ws.Range("D43").Style = "Normal"

$ws.Range("E44").Value = "  +1.53%  "

$ws.Range("E46").Value = "  +0.15%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.5839"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.42%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "129.09"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +3.55%  "

$ws.Range("E49").Value = "  -0.20%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.226"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -3.09%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.07306"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.27%  "
